# Weekly data refresh for "Hortaliza, Vega Modelo de Temuco - Pepino dulce":
# a new price record (week of 2022-02-09) is inserted at the top of the
# historical data block (row 88), pushing every existing record down by
# one row. The sheet's used range grows from A1:R178 to A1:R179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 88..178 down to 89..179, opening up a blank row 88.
$ws.Rows("88:88").Insert()

# Populate the new row 88 with the new weekly record. It carries the same
# market / category / variety / quality / unit / origin metadata as the
# (now shifted-down) row that used to sit at 88, just with a new date and
# new price figures.
$ws.Cells.Item(88, 1).Value2  = 10
$ws.Cells.Item(88, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(88, 3).Value2  = "La Araucanía"
$ws.Cells.Item(88, 4).Value2  = 44601
$ws.Cells.Item(88, 5).Value2  = 9
$ws.Cells.Item(88, 6).Value2  = 100112043
$ws.Cells.Item(88, 7).Value2  = "Pepino dulce"
$ws.Cells.Item(88, 8).Value2  = "Cultivar IV Región"
$ws.Cells.Item(88, 9).Value2  = "Primera"
$ws.Cells.Item(88, 10).Value2 = 35
$ws.Cells.Item(88, 11).Value2 = 19000
$ws.Cells.Item(88, 12).Value2 = 19000
$ws.Cells.Item(88, 13).Value2 = 19000
$ws.Cells.Item(88, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(88, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(88, 16).Value2 = 1056
$ws.Cells.Item(88, 17).Value2 = 18
$ws.Cells.Item(88, 18).Value2 = "Hortaliza"
